# Auto-generated: apply scheduled-runner market-price updates to each Leve sheet.
# Values correspond to updated currentAveragePrice* / LevePrice* / LeveProfit* figures
# refreshed by the scheduled data-sync job (see commit message).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 958.3125
$ws.Range("I12").Value = 95.818184
$ws.Range("K12").Value = 95.818184
$ws.Range("M12").Value = 74.181816
$ws.Range("H18").Value = 724.5
$ws.Range("I18").Value = 724.5
$ws.Range("K18").Value = 724.5
$ws.Range("M18").Value = -440.5
$ws.Range("H19").Value = 785.5714
$ws.Range("I19").Value = 666.3333
$ws.Range("J19").Value = 875
$ws.Range("K19").Value = 666.3333
$ws.Range("L19").Value = 875
$ws.Range("M19").Value = -491.3333
$ws.Range("N19").Value = -1225
$ws.Range("H29").Value = 3259.4167
$ws.Range("J29").Value = 5226.857
$ws.Range("L29").Value = 15680.571
$ws.Range("N29").Value = -16242.571
$ws.Range("H33").Value = 226.3
$ws.Range("I33").Value = 253.125
$ws.Range("J33").Value = 119
$ws.Range("K33").Value = 253.125
$ws.Range("L33").Value = 119
$ws.Range("M33").Value = -24.125
$ws.Range("N33").Value = -577
$ws.Range("H43").Value = 14999.333
$ws.Range("J43").Value = 14999.333
$ws.Range("L43").Value = 14999.333
$ws.Range("N43").Value = -15137.333
$ws.Range("H106").Value = 9699.348
$ws.Range("J106").Value = 14540.667
$ws.Range("L106").Value = 14540.667
$ws.Range("N106").Value = -15802.667
$ws.Range("H113").Value = 7072.143
$ws.Range("J113").Value = 8835.333000000001
$ws.Range("L113").Value = 8835.333000000001
$ws.Range("N113").Value = -15343.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 50003750
$ws.Range("I45").Value = 90910550
$ws.Range("J45").Value = 6548
$ws.Range("K45").Value = 90910550
$ws.Range("L45").Value = 6548
$ws.Range("M45").Value = -90910173
$ws.Range("N45").Value = -7302
$ws.Range("H61").Value = 3718.9614
$ws.Range("I61").Value = 3226.9546
$ws.Range("K61").Value = 3226.9546
$ws.Range("M61").Value = -3014.9546
$ws.Range("H122").Value = 2242.9333
$ws.Range("I122").Value = 1688.8572
$ws.Range("K122").Value = 5066.571599999999
$ws.Range("M122").Value = -2616.571599999999
$ws.Range("H136").Value = 3718.9614
$ws.Range("I136").Value = 3226.9546
$ws.Range("K136").Value = 9680.863799999999
$ws.Range("M136").Value = -7130.863799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2377.6
$ws.Range("I20").Value = 2169.625
$ws.Range("K20").Value = 2169.625
$ws.Range("M20").Value = -1922.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 8450.786
$ws.Range("I69").Value = 5452.5
$ws.Range("K69").Value = 16357.5
$ws.Range("M69").Value = -15546.5
$ws.Range("H72").Value = 8450.786
$ws.Range("I72").Value = 5452.5
$ws.Range("K72").Value = 49072.5
$ws.Range("M72").Value = -45016.5
$ws.Range("H81").Value = 5293.6665
$ws.Range("I81").Value = 1898.75
$ws.Range("J81").Value = 8009.6
$ws.Range("K81").Value = 5696.25
$ws.Range("L81").Value = 24028.8
$ws.Range("M81").Value = -4573.25
$ws.Range("N81").Value = -26274.8
$ws.Range("H84").Value = 5293.6665
$ws.Range("I84").Value = 1898.75
$ws.Range("J84").Value = 8009.6
$ws.Range("K84").Value = 17088.75
$ws.Range("L84").Value = 72086.40000000001
$ws.Range("M84").Value = -11472.75
$ws.Range("N84").Value = -83318.40000000001
$ws.Range("H131").Value = 12154682
$ws.Range("I131").Value = 2426.875
$ws.Range("J131").Value = 24306936
$ws.Range("K131").Value = 7280.625
$ws.Range("L131").Value = 72920808
$ws.Range("M131").Value = -2240.625
$ws.Range("N131").Value = -72930888
$ws.Range("H132").Value = 5716.5835
$ws.Range("J132").Value = 6857.143
$ws.Range("L132").Value = 61714.287
$ws.Range("N132").Value = -66774.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16229.4
$ws.Range("I70").Value = 5407.8335
$ws.Range("J70").Value = 32461.75
$ws.Range("K70").Value = 5407.8335
$ws.Range("L70").Value = 32461.75
$ws.Range("M70").Value = -5137.8335
$ws.Range("N70").Value = -33001.75
$ws.Range("H73").Value = 16229.4
$ws.Range("I73").Value = 5407.8335
$ws.Range("J73").Value = 32461.75
$ws.Range("K73").Value = 5407.8335
$ws.Range("L73").Value = 32461.75
$ws.Range("M73").Value = -4471.8335
$ws.Range("N73").Value = -34333.75
$ws.Range("H122").Value = 3309.0527
$ws.Range("I122").Value = 2704.2666
$ws.Range("K122").Value = 8112.7998
$ws.Range("M122").Value = -5662.7998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4159.9165
$ws.Range("I16").Value = 3292
$ws.Range("J16").Value = 8499.5
$ws.Range("K16").Value = 3292
$ws.Range("L16").Value = 8499.5
$ws.Range("M16").Value = -3122
$ws.Range("N16").Value = -8839.5
$ws.Range("H22").Value = 4077.9092
$ws.Range("I22").Value = 1771
$ws.Range("J22").Value = 6000.3335
$ws.Range("K22").Value = 1771
$ws.Range("L22").Value = 6000.3335
$ws.Range("M22").Value = -1476
$ws.Range("N22").Value = -6590.3335
$ws.Range("H27").Value = 4077.9092
$ws.Range("I27").Value = 1771
$ws.Range("J27").Value = 6000.3335
$ws.Range("K27").Value = 1771
$ws.Range("L27").Value = 6000.3335
$ws.Range("M27").Value = -1664
$ws.Range("N27").Value = -6214.3335
$ws.Range("H46").Value = 4874.75
$ws.Range("J46").Value = 5285.4287
$ws.Range("L46").Value = 5285.4287
$ws.Range("N46").Value = -5661.4287
$ws.Range("H55").Value = 1353006.1
$ws.Range("I55").Value = 2273980.8
$ws.Range("J55").Value = 2243.2
$ws.Range("K55").Value = 2273980.8
$ws.Range("L55").Value = 2243.2
$ws.Range("M55").Value = -2273807.8
$ws.Range("N55").Value = -2589.2
$ws.Range("H68").Value = 4599.294
$ws.Range("I68").Value = 2945.7334
$ws.Range("J68").Value = 17001
$ws.Range("K68").Value = 2945.7334
$ws.Range("L68").Value = 17001
$ws.Range("M68").Value = -2196.7334
$ws.Range("N68").Value = -18499
$ws.Range("H71").Value = 4599.294
$ws.Range("I71").Value = 2945.7334
$ws.Range("J71").Value = 17001
$ws.Range("K71").Value = 14728.667
$ws.Range("L71").Value = 85005
$ws.Range("M71").Value = -10984.667
$ws.Range("N71").Value = -92493
$ws.Range("H132").Value = 4949.619
$ws.Range("J132").Value = 7136.727
$ws.Range("L132").Value = 21410.181
$ws.Range("N132").Value = -26470.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2999.6667
$ws.Range("H65").Value = 2999.6667
